$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose Price (D) / Volume(1h) (E) values were refreshed by the
# "Updated symbol list" GitHub Actions job. Values are stored as literal
# text (e.g. "332.48", "0.96%") matching the source data feed, so force
# each target cell to Text format before writing the new value - this
# prevents Excel from auto-converting the numeric-looking / percent-like
# strings into real numbers.
$targetCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "D4",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "E8",
    "D9",
    "E9",
    "D10",
    "E10",
    "D11",
    "E11",
    "D12",
    "E12",
    "D13",
    "E13",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "E17",
    "D18",
    "E18",
    "D20",
    "E20",
    "D21",
    "E21",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "E26",
    "D27",
    "E27",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D42",
    "E42",
    "D43",
    "E43",
    "D44",
    "E44",
    "E45",
    "D46",
    "E46",
    "E47",
    "D48",
    "E48",
    "E49",
    "D50",
    "E50",
    "E51"
)

foreach ($cellRef in $targetCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "332.48"
$ws.Range("E2").Value = "0.96%"
$ws.Range("D3").Value = "44.02"
$ws.Range("E3").Value = "7.25%"
$ws.Range("D4").Value = "5.793"
$ws.Range("E4").Value = "3.81%"
$ws.Range("D5").Value = "0.08350"
$ws.Range("E5").Value = "2.22%"
$ws.Range("D6").Value = "8.812"
$ws.Range("E6").Value = "0.90%"
$ws.Range("D7").Value = "1.966"
$ws.Range("E7").Value = "-3.89%"
$ws.Range("E8").Value = "-1.88%"
$ws.Range("D9").Value = "0.9326"
$ws.Range("E9").Value = "1.54%"
$ws.Range("D10").Value = "0.1240"
$ws.Range("E10").Value = "-1.75%"
$ws.Range("D11").Value = "0.1954"
$ws.Range("E11").Value = "0.00%"
$ws.Range("D12").Value = "0.09477"
$ws.Range("E12").Value = "0.41%"
$ws.Range("D13").Value = "0.03959"
$ws.Range("E13").Value = "5.78%"
$ws.Range("E14").Value = "0.94%"
$ws.Range("D15").Value = "0.001310"
$ws.Range("E15").Value = "0.36%"
$ws.Range("D16").Value = "0.005922"
$ws.Range("E16").Value = "-4.89%"
$ws.Range("E17").Value = "1.93%"
$ws.Range("D18").Value = "4.505"
$ws.Range("E18").Value = "-0.59%"
$ws.Range("D20").Value = "9.061"
$ws.Range("E20").Value = "8.69%"
$ws.Range("D21").Value = "0.1372"
$ws.Range("E21").Value = "-1.57%"
$ws.Range("D23").Value = "0.04407"
$ws.Range("E23").Value = "-0.10%"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").Value = "-0.03%"
$ws.Range("D25").Value = "0.004350"
$ws.Range("E25").Value = "1.18%"
$ws.Range("E26").Value = "0.75%"
$ws.Range("D27").Value = "0.0003994"
$ws.Range("E27").Value = "0.02%"
$ws.Range("D39").Value = "0.02836"
$ws.Range("E39").Value = "2.78%"
$ws.Range("D40").Value = "0.05821"
$ws.Range("E40").Value = "7.69%"
$ws.Range("D41").Value = "0.007968"
$ws.Range("E41").Value = "4.06%"
$ws.Range("D42").Value = "0.1426"
$ws.Range("E42").Value = "0.87%"
$ws.Range("D43").Value = "0.009087"
$ws.Range("E43").Value = "0.94%"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").Value = "-0.57%"
$ws.Range("E45").Value = "-6.97%"
$ws.Range("D46").Value = "0.00007274"
$ws.Range("E46").Value = "5.50%"
$ws.Range("E47").Value = "-0.08%"
$ws.Range("D48").Value = "0.003335"
$ws.Range("E48").Value = "-6.89%"
$ws.Range("E49").Value = "-0.15%"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").Value = "-0.08%"
$ws.Range("E51").Value = "-0.08%"
